# Add time logs for Sprint 3 (row 9) per contributor.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Durations stored as day-fractions (format [h]:mm):
#   B9 (Zeno Paukner)       -> 25:05
#   C9 (Martin Hausleitner) -> 27:25
#   D9 (Romeo Bhuiyan)      -> 25:18
#   E9 (Jakob Lehner)       -> 18:33
$ws.Range("B9").Value = 1.0451388888888888
$ws.Range("C9").Value = 1.1416666666666666
$ws.Range("D9").Value = 1.0541666666666667
$ws.Range("E9").Value = 0.7729166666666667

$ws.Range("E13").Select()
